$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: extend the merged "Result" header from C1:G1 to C1:J1 ---
# Merging directly over the larger range extends the merge and propagates
# the style of the existing merged cell (s=3) to the newly covered cells.
$ws.Range("C1:J1").Merge()

# --- Row 2: insert the new header cells, copying formatting from existing
#     header cells of the same style before writing the final label text ---

# New cell that should take style s=2 (same as the old G2) -- copy this
# BEFORE G2 itself gets overwritten with the s=1 style below.
$ws.Range("G2").Copy($ws.Range("J2"))

# New cells that should take style s=1 (same as C2/D2/E2/F2 originally)
$ws.Range("C2").Copy($ws.Range("E2"))
$ws.Range("C2").Copy($ws.Range("G2"))
$ws.Range("C2").Copy($ws.Range("H2"))
$ws.Range("C2").Copy($ws.Range("I2"))
$ws.Range("C2").Copy($ws.Range("K2"))

# Now assign the final label text for every header cell in row 2
$ws.Range("D2").Value = "tax_calc_no_glob"
$ws.Range("E2").Value = "solde glob 0"
$ws.Range("F2").Value = "tax_calc_glob_p1"
$ws.Range("G2").Value = "solde p1"
$ws.Range("H2").Value = "tax_calc_glob_p2"
$ws.Range("I2").Value = "solde p2"
$ws.Range("J2").Value = "tax_calc_glob_p1_p2"
$ws.Range("K2").Value = "solde p1_p2"

# --- Row 3: remove the old "v0" value and fill in the new numeric data ---
$ws.Range("B3").ClearContents()

$ws.Range("D3").NumberFormat = "#,##0.00"
$ws.Range("D3").Value = -16207.56
$ws.Range("E3").Value = -1357.34
$ws.Range("F3").NumberFormat = "#,##0.00"
$ws.Range("F3").Value = -16168.27
$ws.Range("G3").Value = -1318.26
$ws.Range("H3").Value = -16207.56
$ws.Range("I3").Value = -1357.34
$ws.Range("J3").Value = -16168.27
$ws.Range("K3").Value = -1318.26

# --- Column widths for the newly introduced columns ---
# (ColumnWidth values chosen so the stored worksheet width lands as close
# as possible to the target 17.1640625 / 20 / 17.33203125 given the
# engine's internal pixel-rounding of column widths)
$ws.Columns("E").ColumnWidth = 16.3
$ws.Columns("G").ColumnWidth = 16.3
$ws.Columns("H").ColumnWidth = 16.3
$ws.Columns("I").ColumnWidth = 16.3
$ws.Columns("J").ColumnWidth = 19.17
$ws.Columns("K").ColumnWidth = 16.45

# --- Update the selected/active cell ---
$ws.Range("I4").Select()
